# Updated cryptos list — refresh Price (D) and Volume(1h) (E) columns,
# plus the RenderToken/THORChain row swap (rows 37-38: B, C, D, E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells are plain text (inlineStr) even when the text looks
# numeric (e.g. "244.80", "36.730.66"). Force the cell to Text format
# before assigning so Excel doesn't auto-coerce the string into a number
# (which would silently drop formatting like trailing zeros).
function Set-Cell($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-Cell "D2" "36.730.66"
Set-Cell "E2" "  +0.79%  "

# Row 3 - Ethereum
Set-Cell "D3" "1.964.86"
Set-Cell "E3" "  +0.90%  "

# Row 4 - TetherUSD
Set-Cell "E4" "  +0.12%  "

# Row 5 - BNB
Set-Cell "D5" "244.80"
Set-Cell "E5" "  +0.27%  "

# Row 6 - XRP
Set-Cell "E6" "  +0.31%  "

# Row 7 - Solana
Set-Cell "D7" "59.25"
Set-Cell "E7" "  +1.35%  "

# Row 8 - USDC
Set-Cell "E8" "  +0.04%  "

# Row 9 - Cardano
Set-Cell "D9" "0.375"
Set-Cell "E9" "  +2.00%  "

# Row 10 - Dogecoin
Set-Cell "D10" "0.0816"
Set-Cell "E10" "  -2.70%  "

# Row 11 - TRON
Set-Cell "E11" "  -0.21%  "

# Row 12 - Avalanche
Set-Cell "D12" "22.38"
Set-Cell "E12" "  +3.17%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-Cell "D13" "2.252.52"
Set-Cell "E13" "  +1.34%  "

# Row 14 - Polygon
Set-Cell "E14" "  +0.41%  "

# Row 15 - Chainlink
Set-Cell "D15" "13.73"
Set-Cell "E15" "  +0.25%  "

# Row 16 - Polkadot
Set-Cell "E16" "  +0.12%  "

# Row 17 - WrappedEther
Set-Cell "D17" "1.960.41"
Set-Cell "E17" "  +0.97%  "

# Row 18 - WrappedBTC
Set-Cell "D18" "36.650.18"
Set-Cell "E18" "  +0.76%  "

# Row 19 - Litecoin
Set-Cell "D19" "70.05"
Set-Cell "E19" "  +0.23%  "

# Row 20 - ShibaInu
Set-Cell "E20" "  -0.77%  "

# Row 21 - Uniswap
Set-Cell "D21" "5.08"
Set-Cell "E21" "  -0.23%  "

# Row 22 - BitcoinCash
Set-Cell "D22" "229.28"
Set-Cell "E22" "  -0.16%  "

# Row 23 - Dai
Set-Cell "E23" "  -0.15%  "

# Row 24 - PancakeSwap
Set-Cell "E24" "  -0.28%  "

# Row 25 - Toncoin
Set-Cell "D25" "2.37"
Set-Cell "E25" "  +3.06%  "

# Row 26 - Kaspa
Set-Cell "D26" "0.146"
Set-Cell "E26" "  +17.69%  "

# Row 27 - Cosmos
Set-Cell "D27" "9.29"
Set-Cell "E27" "  +0.11%  "

# Row 28 - Monero
Set-Cell "D28" "160.74"
Set-Cell "E28" "  -1.24%  "

# Row 29 - EthereumClassic
Set-Cell "E29" "  -0.13%  "

# Row 30 - Stellar
Set-Cell "E30" "  +1.33%  "

# Row 31 - ImmutableX
Set-Cell "D31" "1.16"
Set-Cell "E31" "  -0.29%  "

# Row 32 - Filecoin
Set-Cell "D32" "4.72"
Set-Cell "E32" "  +0.52%  "

# Row 33 - Hedera
Set-Cell "D33" "0.0619"
Set-Cell "E33" "  -1.62%  "

# Row 34 - InternetComputer(DFINITY)
Set-Cell "E34" "  -0.01%  "

# Row 35 - LidoDAOToken
Set-Cell "E35" "  +6.50%  "

# Row 36 - BinanceUSD
Set-Cell "E36" "  +0.17%  "

# Row 37 / 38 - RenderToken <-> THORChain swap (with updated price/volume)
Set-Cell "B37" "THORChain"
Set-Cell "C37" "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-Cell "D37" "6.03"
Set-Cell "E37" "  -3.23%  "

Set-Cell "B38" "RenderToken"
Set-Cell "C38" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-Cell "D38" "3.39"
Set-Cell "E38" "  +11.12%  "

# Row 39 - WEMIXToken
Set-Cell "E39" "  +0.04%  "

# Row 40 - Cronos
Set-Cell "E40" "  +0.05%  "

# Row 41 - HuobiToken
Set-Cell "E41" "  +0.20%  "

# Row 42 - TrustWalletToken
Set-Cell "E42" "  -0.54%  "

# Row 43 - VeChain
Set-Cell "D43" "0.0213"
Set-Cell "E43" "  +1.53%  "

# Row 44 - InjectiveProtocol
Set-Cell "D44" "16.23"
Set-Cell "E44" "  +0.45%  "

# Row 45 - Maker
Set-Cell "D45" "1.359.07"
Set-Cell "E45" "  +0.37%  "

# Row 46 - ARBITRUM
Set-Cell "E46" "  +0.49%  "

# Row 47 - Aave
Set-Cell "D47" "87.94"
Set-Cell "E47" "  -0.09%  "

# Row 48 - FraxShare
Set-Cell "D48" "7.18"
Set-Cell "E48" "  -0.08%  "

# Row 49 - MXToken
Set-Cell "D49" "2.84"
Set-Cell "E49" "  +0.93%  "

# Row 50 - RocketPoolETH
Set-Cell "D50" "2.144.20"
Set-Cell "E50" "  +1.33%  "

# Row 51 - MultiversX
Set-Cell "D51" "43.94"
Set-Cell "E51" "  -3.93%  "
